$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewDNCN")

# Update the pending invoice reference value in D2
$ws.Range("D2").Value = "EADN21-0757"

# Reflect the cell that was last edited/selected
$ws.Range("D2").Select()
